$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26
$ws.Range("B26").Value = 6732773
$ws.Range("F26").Value = 'Suduva Marijampole'
$ws.Range("G26").Value = 'Hegelmann Litauen'
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 1
$ws.Range("L26").Value = 3.8
$ws.Range("M26").Value = 1.533
$ws.Range("N26").Value = 5
$ws.Range("O26").Value = 4.2
$ws.Range("P26").Value = 1.533
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = 1.875
$ws.Range("S26").Value = 1.925
$ws.Range("U26").Value = 1.9
$ws.Range("V26").Value = 1.9
$ws.Range("Y26").Value = 0.5329999999999999
$ws.Range("Z26").Value = 0
$ws.Range("AA26").Value = -0
$ws.Range("AB26").Value = -1
$ws.Range("AC26").Value = 0.8999999999999999

# Row 27
$ws.Range("B27").Value = 6732711
$ws.Range("F27").Value = 'Banga Gargzdai'
$ws.Range("G27").Value = 'FK Zalgiris Vilnius'
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = 4
$ws.Range("L27").Value = 3.6
$ws.Range("M27").Value = 1.571
$ws.Range("N27").Value = 11
$ws.Range("O27").Value = 4.75
$ws.Range("P27").Value = 1.25
$ws.Range("Q27").Value = 1.5
$ws.Range("R27").Value = 1.975
$ws.Range("S27").Value = 1.825
$ws.Range("U27").Value = 1.8
$ws.Range("V27").Value = 2
$ws.Range("Y27").Value = 0.25
$ws.Range("Z27").Value = -1
$ws.Range("AA27").Value = 0.825
$ws.Range("AB27").Value = 0.8
$ws.Range("AC27").Value = -1

# Row 50
$ws.Range("B50").Value = 6732794
$ws.Range("F50").Value = 'FK Siauliai'
$ws.Range("G50").Value = 'FK Dziugas Telsiai'
$ws.Range("H50").Value = 3
$ws.Range("K50").Value = 1.25
$ws.Range("L50").Value = 5
$ws.Range("M50").Value = 9
$ws.Range("N50").Value = 1.25
$ws.Range("O50").Value = 5.25
$ws.Range("P50").Value = 9
$ws.Range("Q50").Value = -1.75
$ws.Range("R50").Value = 2
$ws.Range("S50").Value = 1.8
$ws.Range("T50").Value = 3
$ws.Range("U50").Value = 1.975
$ws.Range("V50").Value = 1.825
$ws.Range("W50").Value = 0.25
$ws.Range("Z50").Value = 1
$ws.Range("AB50").Value = 0
$ws.Range("AC50").Value = -0

# Row 51
$ws.Range("B51").Value = 6732795
$ws.Range("F51").Value = 'Suduva Marijampole'
$ws.Range("G51").Value = 'Banga Gargzdai'
$ws.Range("H51").Value = 1
$ws.Range("K51").Value = 2.15
$ws.Range("L51").Value = 3.2
$ws.Range("M51").Value = 3
$ws.Range("N51").Value = 2.3
$ws.Range("O51").Value = 3.2
$ws.Range("P51").Value = 2.7
$ws.Range("Q51").Value = -0.25
$ws.Range("R51").Value = 2.05
$ws.Range("S51").Value = 1.75
$ws.Range("T51").Value = 2.25
$ws.Range("U51").Value = 1.9
$ws.Range("V51").Value = 1.9
$ws.Range("W51").Value = 1.3
$ws.Range("Z51").Value = 1.05
$ws.Range("AB51").Value = -1
$ws.Range("AC51").Value = 0.8999999999999999

# Row 100
$ws.Range("B100").Value = 6732836
$ws.Range("F100").Value = 'FK Siauliai'
$ws.Range("G100").Value = 'Banga Gargzdai'
$ws.Range("H100").Value = 3
$ws.Range("J100").Value = 'H'
$ws.Range("K100").Value = 1.222
$ws.Range("M100").Value = 9
$ws.Range("N100").Value = 1.363
$ws.Range("P100").Value = 7
$ws.Range("Q100").Value = -1.25
$ws.Range("R100").Value = 1.9
$ws.Range("S100").Value = 1.9
$ws.Range("U100").Value = 1.975
$ws.Range("V100").Value = 1.825
$ws.Range("W100").Value = 0.363
$ws.Range("X100").Value = -1
$ws.Range("Z100").Value = 0.8999999999999999
$ws.Range("AA100").Value = -1
$ws.Range("AB100").Value = 0.9750000000000001
$ws.Range("AC100").Value = -1

# Row 101
$ws.Range("B101").Value = 6732834
$ws.Range("F101").Value = 'Panevezys'
$ws.Range("G101").Value = 'FK Dziugas Telsiai'
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 'D'
$ws.Range("K101").Value = 1.25
$ws.Range("M101").Value = 7.5
$ws.Range("N101").Value = 1.45
$ws.Range("P101").Value = 5
$ws.Range("Q101").Value = -1
$ws.Range("R101").Value = 1.775
$ws.Range("S101").Value = 2.025
$ws.Range("U101").Value = 1.875
$ws.Range("V101").Value = 1.925
$ws.Range("W101").Value = -1
$ws.Range("X101").Value = 3.5
$ws.Range("Z101").Value = -1
$ws.Range("AA101").Value = 1.025
$ws.Range("AB101").Value = -1
$ws.Range("AC101").Value = 0.925

# Row 117
$ws.Range("B117").Value = 7862036
$ws.Range("F117").Value = 'Banga Gargzdai'
$ws.Range("G117").Value = 'FK Zalgiris Vilnius'
$ws.Range("H117").Value = 1
$ws.Range("I117").Value = 4
$ws.Range("J117").Value = 'A'
$ws.Range("K117").Value = 8
$ws.Range("L117").Value = 4.5
$ws.Range("M117").Value = 1.3
$ws.Range("N117").Value = 6.5
$ws.Range("O117").Value = 4.5
$ws.Range("P117").Value = 1.333
$ws.Range("Q117").Value = 1.25
$ws.Range("R117").Value = 2
$ws.Range("S117").Value = 1.8
$ws.Range("U117").Value = 1.825
$ws.Range("V117").Value = 1.975
$ws.Range("X117").Value = -1
$ws.Range("Y117").Value = 0.333
$ws.Range("Z117").Value = -1
$ws.Range("AA117").Value = 0.8
$ws.Range("AB117").Value = 0.825

# Row 118
$ws.Range("B118").Value = 7862911
$ws.Range("F118").Value = 'Hegelmann Litauen'
$ws.Range("G118").Value = 'FK Siauliai'
$ws.Range("H118").Value = 2
$ws.Range("I118").Value = 2
$ws.Range("J118").Value = 'D'
$ws.Range("K118").Value = 2.15
$ws.Range("L118").Value = 3.1
$ws.Range("M118").Value = 3.1
$ws.Range("N118").Value = 2.45
$ws.Range("O118").Value = 2.9
$ws.Range("P118").Value = 3
$ws.Range("Q118").Value = 0
$ws.Range("R118").Value = 1.725
$ws.Range("S118").Value = 2.075
$ws.Range("U118").Value = 2.025
$ws.Range("V118").Value = 1.775
$ws.Range("X118").Value = 1.9
$ws.Range("Y118").Value = -1
$ws.Range("Z118").Value = 0
$ws.Range("AA118").Value = -0
$ws.Range("AB118").Value = 1.025

# Row 140
$ws.Range("B140").Value = 7862924
$ws.Range("E140").Value = 45395.375
$ws.Range("F140").Value = 'FK Transinvest'
$ws.Range("G140").Value = 'Banga Gargzdai'
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 1
$ws.Range("J140").Value = 'A'
$ws.Range("K140").Value = 2.5
$ws.Range("L140").Value = 3.2
$ws.Range("M140").Value = 2.5
$ws.Range("N140").Value = 2.05
$ws.Range("O140").Value = 3.3
$ws.Range("P140").Value = 3.1
$ws.Range("Q140").Value = -0.25
$ws.Range("R140").Value = 1.9
$ws.Range("S140").Value = 1.9
$ws.Range("T140").Value = 2
$ws.Range("U140").Value = 1.775
$ws.Range("V140").Value = 2.025
$ws.Range("W140").Value = -1
$ws.Range("X140").Value = -1
$ws.Range("Y140").Value = 2.1
$ws.Range("Z140").Value = -1
$ws.Range("AA140").Value = 0.8999999999999999
$ws.Range("AB140").Value = -1
$ws.Range("AC140").Value = 1.025

# Row 141
$ws.Range("B141").Value = 7862045
$ws.Range("E141").Value = 45395.45833333334
$ws.Range("F141").Value = 'Suduva Marijampole'
$ws.Range("G141").Value = 'Hegelmann Litauen'
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 1
$ws.Range("J141").Value = 'A'
$ws.Range("K141").Value = 3.4
$ws.Range("M141").Value = 1.95
$ws.Range("N141").Value = 4.5
$ws.Range("O141").Value = 3.8
$ws.Range("P141").Value = 1.571
$ws.Range("Q141").Value = 0.75
$ws.Range("R141").Value = 2
$ws.Range("S141").Value = 1.8
$ws.Range("U141").Value = 1.775
$ws.Range("V141").Value = 2.025
$ws.Range("W141").Value = -1
$ws.Range("X141").Value = -1
$ws.Range("Y141").Value = 0.571
$ws.Range("Z141").Value = -0.5
$ws.Range("AA141").Value = 0.4
$ws.Range("AB141").Value = -1
$ws.Range("AC141").Value = 1.025

# Row 142 is fully removed; subsequent rows shift up implicitly (was already last row)
$ws.Rows("142:142").Delete()
